# Work in progress: fixing not-yet-defined-parent
#
# On the "Concepts" sheet:
#  - Fix wording of the GrainBoundary elucidation/comment cells (row 4).
#  - Insert a new concept row ("SubgrainBoundary") right after GrainBoundary,
#    with a parent ("LowAngleGrainBoundary") that is (deliberately, per the
#    commit message) not yet defined elsewhere in the sheet, to exercise the
#    "undefined parent" test case.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

# --- Fix existing row 4 text ---
$ws.Range("C4").Value = "The boundary of a grain"
$ws.Range("J4").Value = "Test not-yet-defined parent"

# --- Insert new row 5 for the SubgrainBoundary concept ---
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).RowHeight = 13.8

$ws.Range("A5").Value = "SubgrainBoundary"
$ws.Range("C5").Value = "The boundary of a subgrain"
$ws.Range("G5").Value = "LowAngleGrainBoundary"
$ws.Range("J5").Value = "Test of undefined parent"

# --- Update selection to reflect where the author was working ---
$ws.Range("A5").Select() | Out-Null
